# Adds Artisan Commands `showCurve`, `showExtraCurve`, `showEvents`, and
# `showBackgroundEvents` to show/hide curves and events.
#
# In the "Commands" worksheet, insert 4 new rows right before the existing
# "RC Command" section (old row 100) and populate them with the new
# command names (column B) and their descriptions (column C).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Commands")

# Insert 4 blank rows at row 100, pushing the RC Command section (and
# everything below it) down by 4 rows.
$ws.Rows("100:103").Insert()
$ws.Rows("100:103").RowHeight = 13.8

$newCommands = @(
    @("showCurve(<name>,<bool>)", "shows/hides the curve indicated by <name> which is one of { ET, BT, DeltaET, DeltaBT, BackgroundET, BackgroundBT}"),
    @("showExtraCurve(<extra_device>,<curve>,<bool>)", "shows/hides the <curve> (one of {T1,T2}) of the zero-based <extra_device> number"),
    @("showEvents(<event_type>, <bool>)", "shows/hides the events of <event_type> in [1,..,5]"),
    @("showBackgroundEvents(<bool>)", "shows/hides the events of the background profile")
)

$row = 100
foreach ($cmd in $newCommands) {
    $ws.Range("B$row").Value = $cmd[0]
    $ws.Range("C$row").Value = $cmd[1]
    $row = $row + 1
}

# Leave the selection where the author left it after making the edit.
$ws.Activate()
$ws.Range("C101").Select() | Out-Null
